$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 1013: A1013 and D1013 switch from text to numeric values ---
$ws.Range("A1013").Value = 26
$ws.Range("D1013").Value = 2

# --- New rows 1014-1064 ---
# Each entry: row, A value, A isNumeric(1/0), B value(text), C value(text), D value, D isNumeric(1/0)
$rows = @(
    @(1014, "29", 1, "Cloudy", "09/03/2024", "8", 1),
    @(1015, "29", 1, "Cloudy", "09/03/2024", "8", 1),
    @(1016, "29", 1, "Light Rain", "09/03/2024", "8", 1),
    @(1017, "29", 1, "Cloudy", "09/03/2024", "8", 1),
    @(1018, "29", 1, "Light Rain", "09/03/2024", "9", 1),
    @(1019, "30", 1, "Mostly Cloudy", "09/03/2024", "9", 1),
    @(1020, "29", 1, "Light Rain", "09/03/2024", "9", 1),
    @(1021, "30", 1, "Mostly Cloudy", "09/03/2024", "9", 1),
    @(1022, "30", 1, "Mostly Cloudy", "09/03/2024", "9", 1),
    @(1023, "30", 1, "Mostly Cloudy", "09/03/2024", "9", 1),
    @(1024, "30", 1, "Mostly Cloudy", "09/03/2024", "10", 1),
    @(1025, "31", 1, "Mostly Cloudy", "09/03/2024", "10", 1),
    @(1026, "31", 1, "Mostly Cloudy", "09/03/2024", "10", 1),
    @(1027, "31", 1, "Mostly Cloudy", "09/03/2024", "10", 1),
    @(1028, "31", 1, "Mostly Cloudy", "09/03/2024", "10", 1),
    @(1029, "31", 1, "Mostly Cloudy", "09/03/2024", "11", 1),
    @(1030, "31", 1, "Mostly Cloudy", "09/03/2024", "11", 1),
    @(1031, "31", 1, "Mostly Cloudy", "09/03/2024", "11", 1),
    @(1032, "31", 1, "Mostly Cloudy", "09/03/2024", "11", 1),
    @(1033, "32", 1, "Cloudy", "09/03/2024", "11", 1),
    @(1034, "31", 1, "Mostly Cloudy", "09/03/2024", "11", 1),
    @(1035, "32", 1, "Cloudy", "09/03/2024", "12", 1),
    @(1036, "32", 1, "Cloudy", "09/03/2024", "12", 1),
    @(1037, "32", 1, "Cloudy", "09/03/2024", "12", 1),
    @(1038, "32", 1, "Cloudy", "09/03/2024", "12", 1),
    @(1039, "32", 1, "Cloudy", "09/03/2024", "12", 1),
    @(1040, "32", 1, "Cloudy", "09/03/2024", "12", 1),
    @(1041, "32", 1, "Cloudy", "09/03/2024", "13", 1),
    @(1042, "32", 1, "Cloudy", "09/03/2024", "13", 1),
    @(1043, "32", 1, "Cloudy", "09/03/2024", "13", 1),
    @(1044, "32", 1, "Cloudy", "09/03/2024", "13", 1),
    @(1045, "32", 1, "Cloudy", "09/03/2024", "13", 1),
    @(1046, "32", 1, "Cloudy", "09/03/2024", "14", 1),
    @(1047, "32", 1, "Cloudy", "09/03/2024", "14", 1),
    @(1048, "31", 1, "Cloudy", "09/03/2024", "14", 1),
    @(1049, "31", 1, "Cloudy", "09/03/2024", "14", 1),
    @(1050, "31", 1, "Cloudy", "09/03/2024", "14", 1),
    @(1051, "31", 1, "Cloudy", "09/03/2024", "14", 1),
    @(1052, "31", 1, "Cloudy", "09/03/2024", "15", 1),
    @(1053, "31", 1, "Cloudy", "09/03/2024", "15", 1),
    @(1054, "31", 1, "Cloudy", "09/03/2024", "15", 1),
    @(1055, "31", 1, "Cloudy", "09/03/2024", "15", 1),
    @(1056, "31", 1, "Cloudy", "09/03/2024", "15", 1),
    @(1057, "31", 1, "Cloudy", "09/03/2024", "15", 1),
    @(1058, "31", 1, "Cloudy", "09/03/2024", "16", 1),
    @(1059, "31", 1, "Cloudy", "09/03/2024", "16", 1),
    @(1060, "31", 1, "Cloudy", "09/03/2024", "16", 1),
    @(1061, "31", 1, "Cloudy", "09/03/2024", "16", 1),
    @(1062, "31", 1, "Cloudy", "09/03/2024", "16", 1),
    @(1063, "30", 1, "Cloudy", "09/03/2024", "17", 1),
    @(1064, "30", 0, "Cloudy", "09/03/2024", "17", 0)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $aVal = $r[1]
    $aIsNum = $r[2]
    $bVal = $r[3]
    $cVal = $r[4]
    $dVal = $r[5]
    $dIsNum = $r[6]

    $aCell = $ws.Cells.Item($rowNum, 1)
    if ($aIsNum -eq 1) {
        $aCell.Value = [double]$aVal
    } else {
        $aCell.NumberFormat = "@"
        $aCell.Value = $aVal
    }

    # Column B: always plain descriptive text (not numeric-looking), safe to assign directly
    $ws.Cells.Item($rowNum, 2).Value = $bVal

    # Column C: date-like text ("MM/DD/YYYY") must stay literal text, not become a date serial
    $cCell = $ws.Cells.Item($rowNum, 3)
    $cCell.NumberFormat = "@"
    $cCell.Value = $cVal

    $dCell = $ws.Cells.Item($rowNum, 4)
    if ($dIsNum -eq 1) {
        $dCell.Value = [double]$dVal
    } else {
        $dCell.NumberFormat = "@"
        $dCell.Value = $dVal
    }
}
